# Append a freshly-scraped job listing and refresh the "last seen" timestamp
# for every existing row. The new job sorts by priority score to position 9
# (just above the old row 9), so rows 9-14 shift down to 10-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-28 01:29:41"

# ---------------------------------------------------------------------------
# 1) Insert a new row at position 9 - this shifts old rows 9..14 down to 10..15
#    (cell values/styles move with the rows; hyperlink objects do not, so we
#    rebuild the hyperlinks collection from scratch afterwards).
# ---------------------------------------------------------------------------
$ws.Rows(9).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the brand-new row 9 with the newly scraped job.
# ---------------------------------------------------------------------------
$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "急募 【1.6万円/即決】超シンプル仕様の3分タイマーアプリ3本(iOS/Android計6ビルド)の開発"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "10,000 円 ~ 20,000 円 / 募集期間 3 日、取引期間 0 日"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5480298"
$ws.Range("G9").Value = 85
$ws.Range("H9").Value = "◆開発 ◇アプリ"

# ---------------------------------------------------------------------------
# 3) Bump the "取得日時" (fetched-at) timestamp for every data row, 2..15.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# ---------------------------------------------------------------------------
# 4) Rebuild the URL hyperlinks for F2:F15 so they point at the right targets
#    (the engine's row-insert does not itself relocate hyperlink objects).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

for ($r = 2; $r -le 15; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($cell, $cell.Value())
}

# Re-adding hyperlinks creates a fresh (duplicate) "Hyperlink" style entry;
# collapse every URL cell back onto the single shared Hyperlink style so the
# stylesheet doesn't grow an extra near-identical xf.
$ws.Range("F2:F15").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5) Column width tweaks (B: 52 -> 55, D: 30 -> 41). The COM ColumnWidth
#    property is offset from the stored OOXML width by 5/6 of a character,
#    so subtract that to land exactly on the target stored width.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 55 - 5/6
$ws.Columns("D").ColumnWidth = 41 - 5/6
